# 22/10/24 pasos 4 a 7 freecodecamp
#
# Inserts a new "Ahorro paseo" budget line (5%) into the PRESUPUESTO sheet's
# GASTOS VARIABLES block (between "Gastos Extra"/"Libros o Cursos" group and
# "Salud"), which pushes "Salud", "Deudas" and "Navidad" down one row and the
# "TOTAL GASTOS VARIABLES" row from 36 to 37. Also adds a new summary row 38
# with the grand total (AD23+AD37), and minimizes the workbook window.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRESUPUESTO")

# --- Insert a new row 33; everything from old row 33 downward shifts to +1 ---
$ws.Rows.Item(33).Insert()

# --- Populate the new row 33 ("Ahorro paseo", 5%) ---
$ws.Range("A33").Value = "Ahorro paseo"
$ws.Range("B33").Value = 0.05
$ws.Range("AD33").Formula = "=AD`$24*B33"
$ws.Range("AF33").Formula = "=AD33-AE33"

# --- New row 38: grand total = GASTOS FIJOS total (AD23) + GASTOS VARIABLES total (AD37) ---
$ws.Range("AD38").Formula = "=AD23+AD37"

# --- Restore the view state recorded for this sheet ---
$window = $wb.Windows.Item(1)
$window.ScrollColumn = 14
$ws.Range("N1").Select()
$ws.Range("AD38").Select()

# --- Minimize the workbook window (best-effort) ---
$window.WindowState = -4140

$wb.Save()
